# This script inserts a new data row at row 52 on the active sheet,
# shifting the existing rows 52:107 down to 53:108, and populates the
# newly inserted row 52 with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 52 (shifts rows 52-107 down to 53-108)
$ws.Rows(52).Insert()

# Populate the new row 52 with the new record
$ws.Range("A52").Value = 10
$ws.Range("B52").Value = "Vega Modelo de Temuco"
$ws.Range("C52").Value = "La Araucanía"
$ws.Range("D52").Value = 44894
$ws.Range("E52").Value = 9
$ws.Range("F52").Value = "Fruta"
$ws.Range("G52").Value = 100101
$ws.Range("H52").Value = "Berries"
$ws.Range("I52").Value = 100101001
$ws.Range("J52").Value = "Arándano (blue)"
$ws.Range("K52").Value = "Sin especificar"
$ws.Range("L52").Value = "Primera"
$ws.Range("M52").Value = 225
$ws.Range("N52").Value = 2500
$ws.Range("O52").Value = 2600
$ws.Range("P52").Value = 2544
$ws.Range("Q52").Value = "`$/kilo"
$ws.Range("R52").Value = "Región del Maule"
$ws.Range("S52").Value = 2544
$ws.Range("T52").Value = 1
